$wb = $excel.ActiveWorkbook

# --- Conv sheet: update input/kernel values and add formula to E6 ---
$wsConv = $wb.Worksheets.Item("Conv")
$wsConv.Range("B1").Value = 16
$wsConv.Range("E1").Value = 1
$wsConv.Range("B2").Value = 5
$wsConv.Range("E2").Value = 8
$wsConv.Range("E6").Formula = "=E2"

# --- MaxPool sheet: move selection to B1 ---
$wsMaxPool = $wb.Worksheets.Item("MaxPool")
$wsMaxPool.Range("B1").Select()

# --- Conv sheet becomes the active/selected tab with B2 selected ---
$wsConv.Activate()
$wsConv.Range("B2").Select()
